$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (single dot) need to be
# forced to Text so Excel does not silently convert "1.00" -> 1 etc.
# Apply a Text number format, set the value, then clear formatting again so
# the cell keeps matching its original (unstyled) appearance.
$numericTextCells = @(
    @{ Addr = 'D5'; Value = '495.62' }
    @{ Addr = 'D6'; Value = '153.55' }
    @{ Addr = 'D7'; Value = '0.997' }
    @{ Addr = 'D16'; Value = '21.33' }
    @{ Addr = 'D19'; Value = '4.59' }
    @{ Addr = 'D20'; Value = '10.32' }
    @{ Addr = 'D21'; Value = '324.32' }
    @{ Addr = 'D22'; Value = '5.98' }
    @{ Addr = 'D23'; Value = '0.996' }
    @{ Addr = 'D24'; Value = '58.57' }
    @{ Addr = 'D26'; Value = '0.166' }
    @{ Addr = 'D27'; Value = '0.999' }
    @{ Addr = 'D32'; Value = '151.45' }
    @{ Addr = 'D33'; Value = '18.42' }
    @{ Addr = 'D34'; Value = '1.54' }
    @{ Addr = 'D35'; Value = '5.29' }
    @{ Addr = 'D36'; Value = '0.912' }
    @{ Addr = 'D37'; Value = '3.84' }
    @{ Addr = 'D38'; Value = '1.16' }
    @{ Addr = 'D40'; Value = '34.42' }
    @{ Addr = 'D42'; Value = '0.618' }
    @{ Addr = 'D43'; Value = '0.0562' }
    @{ Addr = 'D44'; Value = '0.994' }
    @{ Addr = 'D45'; Value = '4.95' }
    @{ Addr = 'D46'; Value = '269.14' }
    @{ Addr = 'D47'; Value = '0.0951' }
    @{ Addr = 'D48'; Value = '0.0231' }
    @{ Addr = 'D49'; Value = '10.21' }
    @{ Addr = 'D50'; Value = '18.13' }
)
foreach ($item in $numericTextCells) {
    $rng = $ws.Range($item.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.ClearFormats()
}

# Remaining cells (text / URLs / percentages / multi-dot numbers) can be
# assigned directly -- Excel keeps them as text already.
$ws.Range('D2').Value = '57.185.95'
$ws.Range('E2').Value = '  +3.91%  '
$ws.Range('D3').Value = '2.510.94'
$ws.Range('E3').Value = '  +2.23%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E5').Value = '  +3.52%  '
$ws.Range('E6').Value = '  +10.54%  '
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('E8').Value = '  +3.30%  '
$ws.Range('D9').Value = '2.528.43'
$ws.Range('E9').Value = '  +3.10%  '
$ws.Range('E11').Value = '  +5.73%  '
$ws.Range('E12').Value = '  +4.82%  '
$ws.Range('E13').Value = '  +1.61%  '
$ws.Range('D14').Value = '2.949.75'
$ws.Range('E14').Value = '  +2.58%  '
$ws.Range('D15').Value = '57.303.90'
$ws.Range('E15').Value = '  +3.94%  '
$ws.Range('E16').Value = '  +4.41%  '
$ws.Range('E17').Value = '  +3.26%  '
$ws.Range('D18').Value = '2.529.79'
$ws.Range('E18').Value = '  +3.18%  '
$ws.Range('E19').Value = '  +5.87%  '
$ws.Range('E20').Value = '  +4.07%  '
$ws.Range('E21').Value = '  +3.64%  '
$ws.Range('E22').Value = '  +6.23%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('E24').Value = '  +2.34%  '
$ws.Range('E25').Value = '  +2.05%  '
$ws.Range('E26').Value = '  +3.40%  '
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('D28').Value = '2.625.21'
$ws.Range('E28').Value = '  +2.76%  '
$ws.Range('E29').Value = '  +4.56%  '
$ws.Range('D30').Value = '0.0₃0829'
$ws.Range('E30').Value = '  +8.05%  '
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('E32').Value = '  +1.86%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E33').Value = '  +3.09%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E34').Value = '  +4.73%  '
$ws.Range('E35').Value = '  +3.01%  '
$ws.Range('E36').Value = '  +7.53%  '
$ws.Range('E37').Value = '  +7.22%  '
$ws.Range('E38').Value = '  +4.68%  '
$ws.Range('E39').Value = '  +10.59%  '
$ws.Range('E40').Value = '  +3.32%  '
$ws.Range('E41').Value = '  +4.44%  '
$ws.Range('E42').Value = '  +3.04%  '
$ws.Range('E43').Value = '  +4.11%  '
$ws.Range('E44').Value = '  -0.41%  '
$ws.Range('E45').Value = '  +6.92%  '
$ws.Range('E46').Value = '  +5.38%  '
$ws.Range('E47').Value = '  +6.80%  '
$ws.Range('E48').Value = '  +4.68%  '
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('E50').Value = '  +6.58%  '
$ws.Range('D51').Value = '1.904.09'
$ws.Range('E51').Value = '  -1.17%  '
